$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab/workbook entry) to reflect the new "through" date
$ws.Name = "Through 2022-07-26"

# Update the label cell for the July row
$ws.Range("A8").Value = "July (through 07-26)"

# Update the July row (row 8) values
$ws.Range("B8").Value = 33
$ws.Range("D8").Value = 60
$ws.Range("F8").Value = 40
$ws.Range("G8").Value = 116
$ws.Range("H8").Value = 126
$ws.Range("I8").Value = 146

# Update the Total row (row 9) values
$ws.Range("B9").Value = 158
$ws.Range("D9").Value = 450
$ws.Range("F9").Value = 291
$ws.Range("G9").Value = 588
$ws.Range("H9").Value = 886
$ws.Range("I9").Value = 952
